# Auto-generated edit script: updates per commit "chore: update Sheets via scheduled runner"
# Applies numeric corrections to currentAveragePrice / LevePrice / LeveProfit columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5890
$ws.Range("J19").Value = 2180
$ws.Range("L19").Value = 2180
$ws.Range("N19").Value = -2530

$ws.Range("H41").Value = 5370.9165
$ws.Range("I41").Value = 93.333336
$ws.Range("J41").Value = 7130.1113
$ws.Range("K41").Value = 93.333336
$ws.Range("L41").Value = 7130.1113
$ws.Range("M41").Value = 346.666664
$ws.Range("N41").Value = -8010.1113

$ws.Range("H51").Value = 3658.3333
$ws.Range("I51").Value = 4133.3335
$ws.Range("J51").Value = 3500
$ws.Range("K51").Value = 4133.3335
$ws.Range("L51").Value = 3500
$ws.Range("M51").Value = -3649.3335
$ws.Range("N51").Value = -4468

$ws.Range("H53").Value = 45786.316
$ws.Range("I53").Value = 111271.22
$ws.Range("J53").Value = 450.6154
$ws.Range("K53").Value = 111271.22
$ws.Range("L53").Value = 450.6154
$ws.Range("M53").Value = -110634.22
$ws.Range("N53").Value = -1724.6154

$ws.Range("H94").Value = 105
$ws.Range("I94").Value = 105
$ws.Range("K94").Value = 105
$ws.Range("M94").Value = 346

$ws.Range("H96").Value = 417.6842
$ws.Range("I96").Value = 522.9167
$ws.Range("J96").Value = 237.28572
$ws.Range("K96").Value = 1568.7501
$ws.Range("L96").Value = 711.85716
$ws.Range("M96").Value = -195.7501
$ws.Range("N96").Value = -3457.85716

$ws.Range("H101").Value = 44904.668
$ws.Range("I101").Value = 200
$ws.Range("J101").Value = 67257
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 201771
$ws.Range("M101").Value = 1022
$ws.Range("N101").Value = -205015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1379.81
$ws.Range("I32").Value = 1220.9684
$ws.Range("K32").Value = 1220.9684
$ws.Range("M32").Value = -933.9684

$ws.Range("H74").Value = 945.64514
$ws.Range("I74").Value = 798.7037
$ws.Range("J74").Value = 1937.5
$ws.Range("K74").Value = 798.7037
$ws.Range("L74").Value = 1937.5
$ws.Range("M74").Value = 75.29629999999997
$ws.Range("N74").Value = -3685.5

$ws.Range("H77").Value = 945.64514
$ws.Range("I77").Value = 798.7037
$ws.Range("J77").Value = 1937.5
$ws.Range("K77").Value = 3993.5185
$ws.Range("L77").Value = 9687.5
$ws.Range("M77").Value = 374.4814999999999
$ws.Range("N77").Value = -18423.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 430
$ws.Range("I22").Value = 425.33334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 425.33334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -252.33334
$ws.Range("N22").Value = -846

$ws.Range("H99").Value = 1750
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -252
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 41347.383
$ws.Range("I134").Value = 65271.25
$ws.Range("K134").Value = 195813.75
$ws.Range("M134").Value = -193278.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 398.7143
$ws.Range("I22").Value = 358.75
$ws.Range("J22").Value = 452
$ws.Range("K22").Value = 358.75
$ws.Range("L22").Value = 452
$ws.Range("M22").Value = -8.75
$ws.Range("N22").Value = -1152

$ws.Range("H31").Value = 7249488
$ws.Range("I31").Value = 2581.3333
$ws.Range("J31").Value = 15155204
$ws.Range("K31").Value = 2581.3333
$ws.Range("L31").Value = 15155204
$ws.Range("M31").Value = -2286.3333
$ws.Range("N31").Value = -15155794

$ws.Range("H34").Value = 7249488
$ws.Range("I34").Value = 2581.3333
$ws.Range("J34").Value = 15155204
$ws.Range("K34").Value = 2581.3333
$ws.Range("L34").Value = 15155204
$ws.Range("M34").Value = -2379.3333
$ws.Range("N34").Value = -15155608

$ws.Range("H62").Value = 18521106
$ws.Range("I62").Value = 2295.6924
$ws.Range("J62").Value = 66670010
$ws.Range("K62").Value = 2295.6924
$ws.Range("L62").Value = 66670010
$ws.Range("M62").Value = -1671.6924
$ws.Range("N62").Value = -66671258

$ws.Range("H65").Value = 18521106
$ws.Range("I65").Value = 2295.6924
$ws.Range("J65").Value = 66670010
$ws.Range("K65").Value = 11478.462
$ws.Range("L65").Value = 333350050
$ws.Range("M65").Value = -8358.462
$ws.Range("N65").Value = -333356290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 214.3077
$ws.Range("I23").Value = 212.5
$ws.Range("J23").Value = 215.11111
$ws.Range("K23").Value = 637.5
$ws.Range("L23").Value = 645.3333299999999
$ws.Range("M23").Value = -402.5
$ws.Range("N23").Value = -1115.33333

$ws.Range("H113").Value = 534.2245
$ws.Range("I113").Value = 490.2069
$ws.Range("J113").Value = 598.05
$ws.Range("K113").Value = 1470.6207
$ws.Range("L113").Value = 1794.15
$ws.Range("M113").Value = 699.3793000000001
$ws.Range("N113").Value = -6134.15

$ws.Range("H131").Value = 342267.3
$ws.Range("I131").Value = 5265.2173
$ws.Range("J131").Value = 469333.7
$ws.Range("K131").Value = 15795.6519
$ws.Range("L131").Value = 1408001.1
$ws.Range("M131").Value = -10755.6519
$ws.Range("N131").Value = -1418081.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H75").Value = 22876.666
$ws.Range("I75").Value = 10000
$ws.Range("J75").Value = 29315
$ws.Range("K75").Value = 10000
$ws.Range("L75").Value = 29315
$ws.Range("M75").Value = -9064
$ws.Range("N75").Value = -31187

$ws.Range("H78").Value = 22876.666
$ws.Range("I78").Value = 10000
$ws.Range("J78").Value = 29315
$ws.Range("K78").Value = 30000
$ws.Range("L78").Value = 87945
$ws.Range("M78").Value = -25320
$ws.Range("N78").Value = -97305
